# Actualización automática 2025-07-01 08:30:08
#
# Sheet "VENTAS POR GRUPO": several per-category sale values for this
# asesor reset to 0 (and the row-58 "X de 56" counters recomputed to
# "0 de 56" for the affected columns).
#
# Sheet "VENTA MENSUAL": rolling monthly window shifted forward one
# month (marzo/abril/mayo/junio -> abril/mayo/junio/julio); each row's
# C/D/E values take the prior D/E/F values, and F (the new, not-yet-
# reported month) becomes 0. Two totals (E27 and the E58 grand total)
# reflect an additional data correction beyond the pure shift.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO — zero out this period's category sales
# ---------------------------------------------------------------
$ws1ZeroCells = @(
    "Q2","R2",
    "E3",
    "C10","E10","G10","N10",
    "H16","I16","M16","N16",
    "D26","I26","L26","M26","Q26","R26",
    "M27",
    "M29",
    "M31","O31",
    "D34","E34",
    "I40",
    "D44","M44",
    "L45","M45",
    "M50","P50",
    "I51","O51",
    "M55"
)
foreach ($addr in $ws1ZeroCells) {
    $ws1.Range($addr).Value = 0
}

# Row 58 totals ("N de 56") -> "0 de 56" for the columns touched above
$ws1Row58Cells = @("C58","D58","E58","G58","H58","I58","L58","M58","N58","O58","P58","Q58","R58")
foreach ($addr in $ws1Row58Cells) {
    $ws1.Range($addr).Value = "0 de 56"
}

# ---------------------------------------------------------------
# Sheet 2: VENTA MENSUAL — advance the rolling 4-month window
# ---------------------------------------------------------------

# Header: month labels shift left by one month
$ws2.Range("C1").Value = "abril"
$ws2.Range("D1").Value = "mayo"
$ws2.Range("E1").Value = "junio"
$ws2.Range("F1").Value = "julio"

# Column widths adjust along with the new month headers.
# (The engine stores ColumnWidth + 5/6 as the OOXML <col width>, so we
# back that offset out to land on the target stored widths of 14 / 11.)
$ws2.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 10.166666666666666

# Data rows: new C = old D, new D = old E, new E = old F, new F = 0
$ws2.Range("C2").Value = 793.77
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 1191.06
$ws2.Range("F2").Value = 0

$ws2.Range("C3").Value = 2478.81
$ws2.Range("D3").Value = 3316.15
$ws2.Range("E3").Value = 559.83
$ws2.Range("F3").Value = 0

$ws2.Range("C4").Value = 2479.09
$ws2.Range("D4").Value = 0

$ws2.Range("C6").Value = -25.66
$ws2.Range("D6").Value = 0

$ws2.Range("C10").Value = 5689.68
$ws2.Range("D10").Value = 18249.33
$ws2.Range("E10").Value = 3793.28
$ws2.Range("F10").Value = 0

$ws2.Range("C12").Value = 4884.5
$ws2.Range("D12").Value = 1888.52
$ws2.Range("E12").Value = 0

$ws2.Range("C16").Value = 4658.91
$ws2.Range("D16").Value = 6711.78
$ws2.Range("E16").Value = 9019.799999999999
$ws2.Range("F16").Value = 0

$ws2.Range("C26").Value = 0
$ws2.Range("D26").Value = 8756.98
$ws2.Range("E26").Value = 15779.69
$ws2.Range("F26").Value = 0

$ws2.Range("C27").Value = 370.41
$ws2.Range("D27").Value = 3865.18
$ws2.Range("E27").Value = 714.15
$ws2.Range("F27").Value = 0

$ws2.Range("C29").Value = 114.19
$ws2.Range("D29").Value = 787.97
$ws2.Range("E29").Value = 23.76
$ws2.Range("F29").Value = 0

$ws2.Range("E31").Value = 4575.95
$ws2.Range("F31").Value = 0

$ws2.Range("E34").Value = 869.73
$ws2.Range("F34").Value = 0

$ws2.Range("C40").Value = 0
$ws2.Range("E40").Value = 86.40000000000001
$ws2.Range("F40").Value = 0

$ws2.Range("C44").Value = 3838.11
$ws2.Range("D44").Value = 1167.85
$ws2.Range("E44").Value = 1053.12
$ws2.Range("F44").Value = 0

$ws2.Range("C45").Value = 1912.87
$ws2.Range("D45").Value = 737.72
$ws2.Range("E45").Value = 1696.97
$ws2.Range("F45").Value = 0

$ws2.Range("C48").Value = 1709.57
$ws2.Range("D48").Value = 598.58
$ws2.Range("E48").Value = 0

$ws2.Range("D49").Value = 3470.45
$ws2.Range("E49").Value = 0

$ws2.Range("C50").Value = 68.04000000000001
$ws2.Range("D50").Value = 0
$ws2.Range("E50").Value = 84.90000000000001
$ws2.Range("F50").Value = 0

$ws2.Range("C51").Value = 2309.47
$ws2.Range("D51").Value = 0
$ws2.Range("E51").Value = 3819.46
$ws2.Range("F51").Value = 0

$ws2.Range("E55").Value = 25
$ws2.Range("F55").Value = 0

$ws2.Range("C56").Value = 0

$ws2.Range("D57").Value = 808.39
$ws2.Range("E57").Value = 0

# Row 58 grand totals
$ws2.Range("C58").Value = 31281.76
$ws2.Range("D58").Value = 50358.9
$ws2.Range("E58").Value = 43293.1
$ws2.Range("F58").Value = 0
